$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3642.8572
$ws.Range("I64").Value = 3940
$ws.Range("J64").Value = 2900
$ws.Range("K64").Value = 3940
$ws.Range("L64").Value = 2900
$ws.Range("M64").Value = -3692
$ws.Range("N64").Value = -3396
$ws.Range("H67").Value = 3642.8572
$ws.Range("I67").Value = 3940
$ws.Range("J67").Value = 2900
$ws.Range("K67").Value = 3940
$ws.Range("L67").Value = 2900
$ws.Range("M67").Value = -3082
$ws.Range("N67").Value = -4616
$ws.Range("H76").Value = 2658.6956
$ws.Range("I76").Value = 2582.5
$ws.Range("J76").Value = 3166.6667
$ws.Range("K76").Value = 2582.5
$ws.Range("L76").Value = 3166.6667
$ws.Range("M76").Value = -2267.5
$ws.Range("N76").Value = -3796.6667
$ws.Range("H79").Value = 2658.6956
$ws.Range("I79").Value = 2582.5
$ws.Range("J79").Value = 3166.6667
$ws.Range("K79").Value = 2582.5
$ws.Range("L79").Value = 3166.6667
$ws.Range("M79").Value = -1490.5
$ws.Range("N79").Value = -5350.6667
$ws.Range("H112").Value = 1420.89
$ws.Range("I112").Value = 634.1429000000001
$ws.Range("J112").Value = 1480.1075
$ws.Range("K112").Value = 1902.4287
$ws.Range("L112").Value = 4440.3225
$ws.Range("M112").Value = -794.4287000000002
$ws.Range("N112").Value = -6656.3225
$ws.Range("H129").Value = 6757869
$ws.Range("J129").Value = 1023.85297
$ws.Range("L129").Value = 3071.55891
$ws.Range("N129").Value = -13071.55891
$ws.Range("H137").Value = 4172648.2
$ws.Range("I137").Value = 7700192
$ws.Range("J137").Value = 3732.9092
$ws.Range("K137").Value = 23100576
$ws.Range("L137").Value = 11198.7276
$ws.Range("M137").Value = -23098026
$ws.Range("N137").Value = -16298.7276
$ws.Range("H138").Value = 8511.906000000001
$ws.Range("I138").Value = 4707.8887
$ws.Range("J138").Value = 9134.382
$ws.Range("K138").Value = 14123.6661
$ws.Range("L138").Value = 27403.146
$ws.Range("M138").Value = -8983.666100000002
$ws.Range("N138").Value = -37683.146
$ws.Range("H140").Value = 28995
$ws.Range("J140").Value = 28995
$ws.Range("L140").Value = 28995
$ws.Range("N140").Value = -39355

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3827.25
$ws.Range("I63").Value = 2103
$ws.Range("J63").Value = 9000
$ws.Range("K63").Value = 2103
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = -1417
$ws.Range("N63").Value = -10372
$ws.Range("H66").Value = 3827.25
$ws.Range("I66").Value = 2103
$ws.Range("J66").Value = 9000
$ws.Range("K66").Value = 10515
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -7083
$ws.Range("N66").Value = -51864
$ws.Range("H88").Value = 2308.2307
$ws.Range("I88").Value = 1642.8572
$ws.Range("J88").Value = 3084.5
$ws.Range("K88").Value = 1642.8572
$ws.Range("L88").Value = 3084.5
$ws.Range("M88").Value = -1236.8572
$ws.Range("N88").Value = -3896.5
$ws.Range("H91").Value = 2308.2307
$ws.Range("I91").Value = 1642.8572
$ws.Range("J91").Value = 3084.5
$ws.Range("K91").Value = 1642.8572
$ws.Range("L91").Value = 3084.5
$ws.Range("M91").Value = -238.8571999999999
$ws.Range("N91").Value = -5892.5
$ws.Range("H132").Value = 2079.6428
$ws.Range("I132").Value = 1636.4706
$ws.Range("J132").Value = 6600
$ws.Range("K132").Value = 4909.4118
$ws.Range("L132").Value = 19800
$ws.Range("M132").Value = -2379.4118
$ws.Range("N132").Value = -24860

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 41347.832
$ws.Range("I24").Value = 4007.5
$ws.Range("J24").Value = 60018
$ws.Range("K24").Value = 4007.5
$ws.Range("L24").Value = 60018
$ws.Range("M24").Value = -3772.5
$ws.Range("N24").Value = -60488
$ws.Range("H34").Value = 17496.5
$ws.Range("I34").Value = 9980
$ws.Range("K34").Value = 9980
$ws.Range("M34").Value = -9866
$ws.Range("H86").Value = 31503.53
$ws.Range("I86").Value = 1628.5714
$ws.Range("J86").Value = 52416
$ws.Range("K86").Value = 1628.5714
$ws.Range("L86").Value = 52416
$ws.Range("M86").Value = -505.5714
$ws.Range("N86").Value = -54662
$ws.Range("H89").Value = 31503.53
$ws.Range("I89").Value = 1628.5714
$ws.Range("J89").Value = 52416
$ws.Range("K89").Value = 8142.857
$ws.Range("L89").Value = 262080
$ws.Range("M89").Value = -2526.857
$ws.Range("N89").Value = -273312

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2022.2222
$ws.Range("I16").Value = 600
$ws.Range("J16").Value = 3160
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 3160
$ws.Range("M16").Value = -313
$ws.Range("N16").Value = -3734
$ws.Range("H31").Value = 1964224.8
$ws.Range("I31").Value = 2859159.2
$ws.Range("J31").Value = 6555.6875
$ws.Range("K31").Value = 2859159.2
$ws.Range("L31").Value = 6555.6875
$ws.Range("M31").Value = -2858864.2
$ws.Range("N31").Value = -7145.6875
$ws.Range("H34").Value = 1964224.8
$ws.Range("I34").Value = 2859159.2
$ws.Range("J34").Value = 6555.6875
$ws.Range("K34").Value = 2859159.2
$ws.Range("L34").Value = 6555.6875
$ws.Range("M34").Value = -2858957.2
$ws.Range("N34").Value = -6959.6875
$ws.Range("H58").Value = 11114117
$ws.Range("I58").Value = 1584.4482
$ws.Range("K58").Value = 1584.4482
$ws.Range("M58").Value = -1381.4482
$ws.Range("H62").Value = 3405.9375
$ws.Range("J62").Value = 7135.3335
$ws.Range("L62").Value = 7135.3335
$ws.Range("N62").Value = -8383.333500000001
$ws.Range("H65").Value = 3405.9375
$ws.Range("J65").Value = 7135.3335
$ws.Range("L65").Value = 35676.6675
$ws.Range("N65").Value = -41916.6675
$ws.Range("H105").Value = 3352.5454
$ws.Range("I105").Value = 2764.2222
$ws.Range("J105").Value = 6000
$ws.Range("K105").Value = 2764.2222
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = -1017.2222
$ws.Range("N105").Value = -9494
$ws.Range("H113").Value = 2022.2222
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 3160
$ws.Range("K113").Value = 600
$ws.Range("L113").Value = 3160
$ws.Range("M113").Value = 1570
$ws.Range("N113").Value = -7500
$ws.Range("H132").Value = 3881.1765
$ws.Range("I132").Value = 4746.2856
$ws.Range("J132").Value = 3275.6
$ws.Range("K132").Value = 14238.8568
$ws.Range("L132").Value = 9826.799999999999
$ws.Range("M132").Value = -11708.8568
$ws.Range("N132").Value = -14886.8
$ws.Range("H136").Value = 11114117
$ws.Range("I136").Value = 1584.4482
$ws.Range("K136").Value = 4753.3446
$ws.Range("M136").Value = -2203.3446
$ws.Range("H141").Value = 31806.156
$ws.Range("J141").Value = 32348.291
$ws.Range("L141").Value = 32348.291
$ws.Range("N141").Value = -42708.291

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 333850.16
$ws.Range("I4").Value = 800320
$ws.Range("J4").Value = 657.4286
$ws.Range("K4").Value = 2400960
$ws.Range("L4").Value = 1972.2858
$ws.Range("M4").Value = -2400848
$ws.Range("N4").Value = -2196.2858
$ws.Range("H107").Value = 1426.0625
$ws.Range("I107").Value = 488.6
$ws.Range("J107").Value = 1852.1818
$ws.Range("K107").Value = 1465.8
$ws.Range("L107").Value = 5556.5454
$ws.Range("M107").Value = 454.1999999999998
$ws.Range("N107").Value = -9396.545399999999
$ws.Range("H113").Value = 4348950
$ws.Range("I113").Value = 8334149.5
$ws.Range("J113").Value = 1459.091
$ws.Range("K113").Value = 25002448.5
$ws.Range("L113").Value = 4377.272999999999
$ws.Range("M113").Value = -25000278.5
$ws.Range("N113").Value = -8717.272999999999
$ws.Range("H129").Value = 20679.852
$ws.Range("I129").Value = 2746.7273
$ws.Range("J129").Value = 33008.875
$ws.Range("K129").Value = 8240.1819
$ws.Range("L129").Value = 99026.625
$ws.Range("M129").Value = -3240.1819
$ws.Range("N129").Value = -109026.625
$ws.Range("H131").Value = 1485.2549
$ws.Range("J131").Value = 1315.5526
$ws.Range("L131").Value = 3946.6578
$ws.Range("N131").Value = -14026.6578

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 60004
$ws.Range("J4").Value = 60004
$ws.Range("L4").Value = 60004
$ws.Range("N4").Value = -60228
$ws.Range("H70").Value = 3966.5908
$ws.Range("I70").Value = 3960.3125
$ws.Range("J70").Value = 3983.3333
$ws.Range("K70").Value = 3960.3125
$ws.Range("L70").Value = 3983.3333
$ws.Range("M70").Value = -3690.3125
$ws.Range("N70").Value = -4523.3333
$ws.Range("H73").Value = 3966.5908
$ws.Range("I73").Value = 3960.3125
$ws.Range("J73").Value = 3983.3333
$ws.Range("K73").Value = 3960.3125
$ws.Range("L73").Value = 3983.3333
$ws.Range("M73").Value = -3024.3125
$ws.Range("N73").Value = -5855.3333
$ws.Range("H80").Value = 3355.7727
$ws.Range("I80").Value = 3264.5789
$ws.Range("J80").Value = 3933.3333
$ws.Range("K80").Value = 3264.5789
$ws.Range("L80").Value = 3933.3333
$ws.Range("M80").Value = -2266.5789
$ws.Range("N80").Value = -5929.3333
$ws.Range("H83").Value = 3355.7727
$ws.Range("I83").Value = 3264.5789
$ws.Range("J83").Value = 3933.3333
$ws.Range("K83").Value = 16322.8945
$ws.Range("L83").Value = 19666.6665
$ws.Range("M83").Value = -11330.8945
$ws.Range("N83").Value = -29650.6665

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 19359.2
$ws.Range("J59").Value = 19359.2
$ws.Range("L59").Value = 19359.2
$ws.Range("N59").Value = -20667.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 22691.666
$ws.Range("J86").Value = 22691.666
$ws.Range("L86").Value = 22691.666
$ws.Range("N86").Value = -24937.666
$ws.Range("H89").Value = 22691.666
$ws.Range("J89").Value = 22691.666
$ws.Range("L89").Value = 113458.33
$ws.Range("N89").Value = -124690.33
$ws.Range("H101").Value = 23943.143
$ws.Range("J101").Value = 23943.143
$ws.Range("L101").Value = 23943.143
$ws.Range("N101").Value = -30433.143
$ws.Range("H132").Value = 4371948.5
$ws.Range("I132").Value = 6674467.5
$ws.Range("J132").Value = 54725.375
$ws.Range("K132").Value = 20023402.5
$ws.Range("L132").Value = 164176.125
$ws.Range("M132").Value = -20020872.5
$ws.Range("N132").Value = -169236.125
